$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New column D: custom width (~18.57 chars, matches the narrow date column the
# author added alongside the new "data type" column)
$ws.Columns.Item(4).ColumnWidth = 17.67

# Date/time value used for each of the new D cells (serial 42016.51458333333
# == 2015-01-12 12:21:00), formatted with the built-in short-date/time format.
$dateValue = 42016.51458333333
$rows = @(4, 5, 6, 8, 9, 10, 12, 13, 14)
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 4)
    $cell.NumberFormat = "m/d/yy h:mm"
    $cell.Value = $dateValue
}

# Move the active selection to D14, matching the final cursor position
$ws.Range("D14").Select()
